# Updates cryptocurrency price/volume data in the "cryptos" sheet,
# matching a scrape refresh commit (GitHub Actions, 2023-05-20).
# Each Price (D) / Volume(1h) (E) cell is rewritten as text, using a
# leading apostrophe for values that would otherwise be auto-parsed as
# numbers (e.g. "1.015"), so the stored cell stays textual exactly like
# the source data (and preserves values such as "0.000008859" that Excel
# would otherwise render in scientific notation as a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.363.10"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "1.840.58"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +1.39%  "
$ws.Range("D5").Value = "'315.07"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "'0.07469"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.8858"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "'20.51"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.846.62"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "'0.07378"
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").Value = "'5.488"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "'93.31"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "'6.578"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'0.000008859"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "27.388.27"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "'10.72"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").Value = "2.070.01"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'1.911"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'152.01"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "'18.66"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'2.182"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'5.282"
$ws.Range("D30").Value = "'117.99"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "'0.7622"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "'1.180"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "'4.564"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "'2.945"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").Value = "'1.105"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").Value = "'0.05373"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").Value = "'0.01962"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'3.000"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").Value = "'0.5357"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'2.380"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").Value = "'0.1668"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'8.553"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "'0.4984"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "'10.51"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'1.013"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "'105.19"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'1.681"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "'0.06332"
$ws.Range("E51").Value = "  +0.75%  "
